$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 150. This shifts the existing rows
# 150-153 down to 151-154 (dimension grows from A1:R153 to A1:R154),
# matching the diff which pushes the old rows down and adds a brand
# new record in what becomes row 150.
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with its data.
$ws.Cells.Item(150, 1).Value = 7
$ws.Cells.Item(150, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(150, 3).Value = "Ñuble"
$ws.Cells.Item(150, 4).Value = 44448
$ws.Cells.Item(150, 5).Value = 16
$ws.Cells.Item(150, 6).Value = 100112023
$ws.Cells.Item(150, 7).Value = "Brócoli"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 360
$ws.Cells.Item(150, 11).Value = 700
$ws.Cells.Item(150, 12).Value = 750
$ws.Cells.Item(150, 13).Value = 725
$ws.Cells.Item(150, 14).Value = "$/unidad"
$ws.Cells.Item(150, 15).Value = "Región del Maule"
$ws.Cells.Item(150, 16).Value = 725
$ws.Cells.Item(150, 17).Value = 1
$ws.Cells.Item(150, 18).Value = "Hortaliza"

# Give the date cell the same number format/style as the rest of
# column D (style index 2, custom date format) so it matches its
# siblings.
$ws.Cells.Item(150, 4).NumberFormat = $ws.Cells.Item(151, 4).NumberFormat
